$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sprint1")
$ws2 = $wb.Worksheets.Item("Sprint2")

# --- Sprint2: bring formatting for column A (and full row pattern) down into rows 38-48 ---
# Row 37 already has the correct per-column style pattern (A=28,B=28,C=28,D=29,E=30,F=29,G=30,H=30,I=29)
$ws2.Range("A37:I37").Copy()
$ws2.Range("A38:A48").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("C37:I37").Copy()
$ws2.Range("C44:I48").PasteSpecial(-4122)   # xlPasteFormats

# Fill in the "Numero d'activitat" column (A) for rows 38-48
$ws2.Range("A38").Value = 23
$ws2.Range("A39").Value = 24
$ws2.Range("A40").Value = 25
$ws2.Range("A41").Value = 26
$ws2.Range("A42").Value = 27
$ws2.Range("A43").Value = 28
$ws2.Range("A44").Value = 29
$ws2.Range("A45").Value = 30
$ws2.Range("A46").Value = 31
$ws2.Range("A47").Value = 32
$ws2.Range("A48").Value = 33

# New rows 42 and 43 (new tasks for Marc Martin / Leonard Craciun)
$ws2.Range("B42").Value = 2
$ws2.Range("C42").Value = "Marc Martin"
$ws2.Range("D42").Value = "marcmartin60"
$ws2.Range("E42").Value = "Vista de creacion de obras"
$ws2.Range("F42").Value = "Leonard Craciun"
$ws2.Range("G42").Value = 2
$ws2.Range("H42").Value = 1
$ws2.Range("I42").Value = 1.5

$ws2.Range("B43").Value = 2
$ws2.Range("C43").Value = "Marc Martin"
$ws2.Range("D43").Value = "marcmartin60"
$ws2.Range("E43").Value = "Vista de modificacion de obras"
$ws2.Range("F43").Value = "Leonard Craciun"
$ws2.Range("G43").Value = 1
$ws2.Range("H43").Value = 1
$ws2.Range("I43").Value = 1

# Rows 44-48: only activity number (A) and Sprint (B) are filled, rest left blank
$ws2.Range("B44").Value = 2
$ws2.Range("B45").Value = 2
$ws2.Range("B46").Value = 2
$ws2.Range("B47").Value = 2
$ws2.Range("B48").Value = 2

# Update the active selection / scroll position left after editing, as in the source workbook
$null = $ws2.Range("C44").Select()
